# Bug fix: empty footnotes (notes that carry no real annotation text --
# just the auto-numbered footnote reference mark itself) should not be
# left dangling in the document. Walk the footnote collection back to
# front (so removing one doesn't disturb the indices of the ones still
# to be checked) and delete any footnote whose content is empty.

$d = $word.ActiveDocument
$footnotes = $d.Footnotes

for ($i = $footnotes.Count; $i -ge 1; $i--) {
    $note = $footnotes.Item($i)
    $text = $note.Range.Text

    # The only thing present in an "empty" footnote is the automatic
    # reference-mark character itself, so a trimmed length of 0 or 1
    # means there is no real note content.
    if ($text.Trim().Length -le 1) {
        $note.Delete()
    }
}
